# Add the four new "get ..." API request rows to Sheet1 (A11:B14),
# matching the yellow/green/blue-fill style already used by the rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 12-14 need the same fill/style formatting that row 11 already has
# (A11/B11 exist as empty, pre-styled cells in the original sheet).
# Copy that formatting down before filling in the values.
$ws.Range("A11:B11").Copy() | Out-Null
$ws.Range("A12:B14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A11").Value = "get saved videos"
$ws.Range("B11").Value = 110

$ws.Range("A12").Value = "get saved playlists"
$ws.Range("B12").Value = 111

$ws.Range("A13").Value = "get videos by random category"
$ws.Range("B13").Value = 112

$ws.Range("A14").Value = "get video history "
$ws.Range("B14").Value = 113

# Mirror the author's final cursor position after entering the new rows.
$ws.Range("B15").Select() | Out-Null
